$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 74, shifting the existing
# rows 74:93 down to 77:96 (new weekly price update pushed to the top
# of this date-ordered block).
$ws.Rows("74:76").Insert()

# Row 74 - Especial / Provincia de Limarí / 2021-11-05
$ws.Range("A74").Value = 8
$ws.Range("B74").Value = "Terminal La Palmera de La Serena"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44505
$ws.Range("E74").Value = 4
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100107
$ws.Range("H74").Value = "Otros"
$ws.Range("I74").Value = 100107002
$ws.Range("J74").Value = "Chirimoya"
$ws.Range("K74").Value = "Cultivar IV Región"
$ws.Range("L74").Value = "Especial"
$ws.Range("M74").Value = 300
$ws.Range("N74").Value = 2000
$ws.Range("O74").Value = 2100
$ws.Range("P74").Value = 2050
$ws.Range("Q74").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R74").Value = "Provincia de Limarí"
$ws.Range("S74").Value = 2050
$ws.Range("T74").Value = 1

# Row 75 - Primera / Provincia de Limarí / 2021-11-05
$ws.Range("A75").Value = 8
$ws.Range("B75").Value = "Terminal La Palmera de La Serena"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44505
$ws.Range("E75").Value = 4
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100107
$ws.Range("H75").Value = "Otros"
$ws.Range("I75").Value = 100107002
$ws.Range("J75").Value = "Chirimoya"
$ws.Range("K75").Value = "Cultivar IV Región"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 300
$ws.Range("N75").Value = 1700
$ws.Range("O75").Value = 1800
$ws.Range("P75").Value = 1750
$ws.Range("Q75").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R75").Value = "Provincia de Limarí"
$ws.Range("S75").Value = 1750
$ws.Range("T75").Value = 1

# Row 76 - Segunda / Provincia de Limarí / 2021-11-05
$ws.Range("A76").Value = 8
$ws.Range("B76").Value = "Terminal La Palmera de La Serena"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44505
$ws.Range("E76").Value = 4
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100107
$ws.Range("H76").Value = "Otros"
$ws.Range("I76").Value = 100107002
$ws.Range("J76").Value = "Chirimoya"
$ws.Range("K76").Value = "Cultivar IV Región"
$ws.Range("L76").Value = "Segunda"
$ws.Range("M76").Value = 300
$ws.Range("N76").Value = 1200
$ws.Range("O76").Value = 1300
$ws.Range("P76").Value = 1250
$ws.Range("Q76").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R76").Value = "Provincia de Limarí"
$ws.Range("S76").Value = 1250
$ws.Range("T76").Value = 1
